$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 4 - duplicate of row 2 (Kings XI Punjab match)
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 24 2020"
$ws.Range("C4").Value = "Kings XI won by 97 runs"
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = "AB de Villiers "
Set-TextValue $ws.Range("G4") "28"
Set-TextValue $ws.Range("H4") "18"
Set-TextValue $ws.Range("I4") "4"
Set-TextValue $ws.Range("J4") "1"
Set-TextValue $ws.Range("K4") "155.55"

# Row 5 - duplicate of row 3 (Sunrisers Hyderabad match)
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " September 21 2020"
$ws.Range("C5").Value = "RCB won by 10 runs"
$ws.Range("D5").Value = "Royal Challengers Bangalore"
$ws.Range("E5").Value = "Sunrisers Hyderabad"
$ws.Range("F5").Value = "AB de Villiers "
Set-TextValue $ws.Range("G5") "51"
Set-TextValue $ws.Range("H5") "30"
Set-TextValue $ws.Range("I5") "4"
Set-TextValue $ws.Range("J5") "2"
Set-TextValue $ws.Range("K5") "170.00"
